# Dev count run 20251205
# Updates the SCM Report Summary workbook:
#   - Sheet "SCM Report Summary": refresh the report-generation time and
#     bump the GitHub contributor totals (a new committer was detected).
#   - Sheet "GitHub Details": add a row for the newly-seen committer
#     "Veracode Dev Count" on the existing repo, and add a row for the
#     "julian-veracode/GitHubApp-with-config" repo/committer pairing.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("SCM Report Summary")
$wsDetails = $wb.Worksheets.Item("GitHub Details")

# --- SCM Report Summary sheet -------------------------------------------
# Time of report
$wsSummary.Range("B3").Value = "10:05:53 AM"

# Total unique contributors across GitHub
$wsSummary.Range("B5").Value = 2

# Total unique across All SCM Platforms
$wsSummary.Range("B7").Value = 2

# --- GitHub Details sheet -------------------------------------------------
# Insert a new row 3 for the extra committer found on the first repo.
$wsDetails.Rows.Item(3).Insert()
$wsDetails.Range("A3").Value = "julz0815/dev-count-test-headless"
$wsDetails.Range("B3").Value = "Veracode Dev Count"
$wsDetails.Range("C3").Value = "dev-count@veracode.com"

# Append the newly observed repo/committer row at the end of the table.
$wsDetails.Range("A11").Value = "julian-veracode/GitHubApp-with-config"
$wsDetails.Range("B11").Value = "Julian Totzek-Hallhuber"
$wsDetails.Range("C11").Value = "j.totzek@gmail.com"
